$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new BOM row for the reverse-bias protection Schottky diode (D1)
# right above the IC1 row, pushing every row below it down by one.
# First, shift the existing cell formatting of rows 7-13 down to rows 8-14
# (copy-format only, so the existing style table is reused unchanged).
$ws.Range("A7:F13").Copy()
$ws.Range("A8:F14").PasteSpecial(-4122)

# Row 7: new Schottky diode entry
$ws.Cells.Item(7,1).Value2 = "Schottky Diode"
$ws.Cells.Item(7,2).Value2 = "RB715UMTL"
$ws.Cells.Item(7,3).Value2 = "D1"
$ws.Cells.Item(7,4).Value2 = ""
$ws.Cells.Item(7,5).Value2 = 1
$ws.Cells.Item(7,6).Value2 = "DAN217UMTL"

# Row 8: IC1 (was row 7)
$ws.Cells.Item(8,1).Value2 = "Integrated Circuit"
$ws.Cells.Item(8,2).Value2 = "STM32G031J6M6"
$ws.Cells.Item(8,3).Value2 = "IC1"
$ws.Cells.Item(8,4).Value2 = ""
$ws.Cells.Item(8,5).Value2 = 1
$ws.Cells.Item(8,6).Value2 = "SOIC127P600X175-8N"

# Row 9: IC2 / LDO (was row 8)
$ws.Cells.Item(9,1).Value2 = "LDO"
$ws.Cells.Item(9,2).Value2 = "MIC5366-3.3YC5-TR"
$ws.Cells.Item(9,3).Value2 = "IC2"
$ws.Cells.Item(9,4).Value2 = "3.3V/150mA"
$ws.Cells.Item(9,5).Value2 = 1
$ws.Cells.Item(9,6).Value2 = "SOT65P210X110-5N"

# Row 10: IC3 (was row 9)
$ws.Cells.Item(10,1).Value2 = "Integrated Circuit"
$ws.Cells.Item(10,2).Value2 = "TB67H450FNG,EL"
$ws.Cells.Item(10,3).Value2 = "IC3"
$ws.Cells.Item(10,4).Value2 = ""
$ws.Cells.Item(10,5).Value2 = 1
$ws.Cells.Item(10,6).Value2 = "SOIC127P600X175-9N"

# Row 11: J2 connector (was row 10)
$ws.Cells.Item(11,1).Value2 = "Connector"
$ws.Cells.Item(11,2).Value2 = "B2B-PH-K-S_LF__SN_"
$ws.Cells.Item(11,3).Value2 = "J2"
$ws.Cells.Item(11,4).Value2 = ""
$ws.Cells.Item(11,5).Value2 = 1
$ws.Cells.Item(11,6).Value2 = "SHDR2W50P0X200_1X2_590X450X600P"

# Row 12: Switch S1 (was row 11)
$ws.Cells.Item(12,1).Value2 = "Switch"
$ws.Cells.Item(12,2).Value2 = "PTS815_SJM_250_SMTR_LFS"
$ws.Cells.Item(12,3).Value2 = "S1"
$ws.Cells.Item(12,4).Value2 = ""
$ws.Cells.Item(12,5).Value2 = 1
$ws.Cells.Item(12,6).Value2 = "EVP-BT3G4A000"

# Row 13: Jumpers (was row 12)
$ws.Cells.Item(13,1).Value2 = "Jumper (0 ohm)"
$ws.Cells.Item(13,2).Value2 = "RC0805FR-070RL"
$ws.Cells.Item(13,3).Value2 = "SB1, SB2, SB4, SB5, SB6"
$ws.Cells.Item(13,4).Value2 = ""
$ws.Cells.Item(13,5).Value2 = 5
$ws.Cells.Item(13,6).Value2 = "RESC2012X60N"

# Row 14: Crystal / Y1 (was row 13)
$ws.Cells.Item(14,1).Value2 = "Crystal or Oscillator"
$ws.Cells.Item(14,2).Value2 = "ASEKDV-32.768kHz-LC-T"
$ws.Cells.Item(14,3).Value2 = "Y1"
$ws.Cells.Item(14,4).Value2 = "32.768kHz"
$ws.Cells.Item(14,5).Value2 = 1
$ws.Cells.Item(14,6).Value2 = "ASEKDV-32.768kHz-LC-T"
